$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.358.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.13%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3712"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.06"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.83%  "

$ws.Range("E9").Value = "  -3.09%  "

$ws.Range("E10").Value = "  -4.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.049"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.919"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.565.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.21%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.238"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5284"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.00%  "

$ws.Range("E24").Value = "  -3.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.351.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.387"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.820"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.66%  "

$ws.Range("E28").Value = "  -4.25%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.14%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.733.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.200"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "

$ws.Range("E35").Value = "  -5.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08434"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.99%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2323"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.523"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06405"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.288"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5974"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.760"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.094"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.267"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
